# Insert a new weekly price-report row at row 33 (pushing the existing
# rows 33:90 down to 34:91), then populate the new row with the latest
# "Poroto verde" (Vega Modelo de Temuco) observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 33:90 down to 34:91 by inserting a blank row at 33.
$ws.Rows("33:33").Insert()

# Fill in the newly inserted row 33 with the new weekly record.
$ws.Range("A33").Value2 = 10
$ws.Range("B33").Value2 = "Vega Modelo de Temuco"
$ws.Range("C33").Value2 = "La Araucanía"
$ws.Range("D33").Value2 = 44519
$ws.Range("E33").Value2 = 9
$ws.Range("F33").Value2 = 100112031
$ws.Range("G33").Value2 = "Poroto verde"
$ws.Range("H33").Value2 = "Sin especificar"
$ws.Range("I33").Value2 = "Primera"
$ws.Range("J33").Value2 = 300
$ws.Range("K33").Value2 = 3000
$ws.Range("L33").Value2 = 3000
$ws.Range("M33").Value2 = 3000
$ws.Range("N33").Value2 = '$/kilo'
$ws.Range("O33").Value2 = "Región Metropolitana"
$ws.Range("P33").Value2 = 3000
$ws.Range("Q33").Value2 = 1
$ws.Range("R33").Value2 = "Hortaliza"
